$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1084.7142
$ws.Range("I2").Value = 918.6
$ws.Range("K2").Value = 918.6
$ws.Range("M2").Value = -805.6
$ws.Range("H33").Value = 125.63636
$ws.Range("I33").Value = 122.75
$ws.Range("J33").Value = 133.33333
$ws.Range("K33").Value = 122.75
$ws.Range("L33").Value = 133.33333
$ws.Range("M33").Value = 106.25
$ws.Range("N33").Value = -591.3333299999999
$ws.Range("H41").Value = 335.8
$ws.Range("I41").Value = 335.8
$ws.Range("K41").Value = 335.8
$ws.Range("M41").Value = 104.2
$ws.Range("H86").Value = 6656.3335
$ws.Range("I86").Value = 6656.3335
$ws.Range("K86").Value = 6656.3335
$ws.Range("M86").Value = -5533.3335
$ws.Range("H89").Value = 6656.3335
$ws.Range("I89").Value = 6656.3335
$ws.Range("K89").Value = 33281.6675
$ws.Range("M89").Value = -27665.6675
$ws.Range("H116").Value = 4899
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H118").Value = 173
$ws.Range("I118").Value = 173
$ws.Range("K118").Value = 519
$ws.Range("M118").Value = 1138
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H137").Value = 10820.182
$ws.Range("I137").Value = 10431.286
$ws.Range("J137").Value = 11500.75
$ws.Range("K137").Value = 31293.858
$ws.Range("L137").Value = 34502.25
$ws.Range("M137").Value = -28743.858
$ws.Range("N137").Value = -39602.25
$ws.Range("H141").Value = 1200
$ws.Range("I141").Value = 1200
$ws.Range("K141").Value = 3600
$ws.Range("M141").Value = 1580

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 468.7143
$ws.Range("I32").Value = 468.7143
$ws.Range("K32").Value = 468.7143
$ws.Range("M32").Value = -181.7143
$ws.Range("H45").Value = 12499
$ws.Range("I45").Value = 12499
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 12499
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H63").Value = 4605.4
$ws.Range("I63").Value = 4605.4
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4605.4
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 4605.4
$ws.Range("I66").Value = 4605.4
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 23027
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 5536.8335
$ws.Range("J74").Value = 5886.364
$ws.Range("L74").Value = 5886.364
$ws.Range("N74").Value = -7634.364
$ws.Range("H77").Value = 5536.8335
$ws.Range("J77").Value = 5886.364
$ws.Range("L77").Value = 29431.82
$ws.Range("N77").Value = -38167.82
$ws.Range("H88").Value = 2423.5
$ws.Range("H91").Value = 2423.5
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1000
$ws.Range("L86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("M86").Value = 123
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5000
$ws.Range("L89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("M89").Value = 616
$ws.Range("H109").Value = 60000
$ws.Range("J109").Value = 60000
$ws.Range("L109").Value = 60000
$ws.Range("N109").Value = -62774

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1118.6364
$ws.Range("I22").Value = 1162.875
$ws.Range("J22").Value = 1000.6667
$ws.Range("K22").Value = 1162.875
$ws.Range("L22").Value = 1000.6667
$ws.Range("M22").Value = -812.875
$ws.Range("N22").Value = -1700.6667
$ws.Range("H31").Value = 6691.6665
$ws.Range("I31").Value = 4814.2856
$ws.Range("K31").Value = 4814.2856
$ws.Range("M31").Value = -4519.2856
$ws.Range("H34").Value = 6691.6665
$ws.Range("I34").Value = 4814.2856
$ws.Range("K34").Value = 4814.2856
$ws.Range("M34").Value = -4612.2856
$ws.Range("H58").Value = 6222.1113
$ws.Range("J58").Value = 6499.875
$ws.Range("L58").Value = 6499.875
$ws.Range("N58").Value = -6905.875
$ws.Range("H86").Value = 15500
$ws.Range("I86").Value = 15500
$ws.Range("K86").Value = 15500
$ws.Range("M86").Value = -14377
$ws.Range("H89").Value = 15500
$ws.Range("I89").Value = 15500
$ws.Range("K89").Value = 77500
$ws.Range("M89").Value = -71884
$ws.Range("H134").Value = 5631.2
$ws.Range("I134").Value = 2044.5714
$ws.Range("K134").Value = 6133.7142
$ws.Range("M134").Value = -3598.7142
$ws.Range("H136").Value = 6222.1113
$ws.Range("J136").Value = 6499.875
$ws.Range("L136").Value = 19499.625
$ws.Range("N136").Value = -24599.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -12746
$ws.Range("H84").Value = 3500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -42732

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 26632.666
$ws.Range("J98").Value = 26632.666
$ws.Range("L98").Value = 26632.666
$ws.Range("N98").Value = -32622.666
$ws.Range("H102").Value = 3418.4285
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H107").Value = 1524.75
$ws.Range("I107").Value = 1599.7142
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1599.7142
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 320.2858000000001
$ws.Range("N107").Value = -4840
$ws.Range("H113").Value = 12498
$ws.Range("I113").Value = 12498
$ws.Range("K113").Value = 12498
$ws.Range("M113").Value = -10328

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9998
$ws.Range("I16").Value = 9998
$ws.Range("K16").Value = 9998
$ws.Range("M16").Value = -9828
$ws.Range("H40").Value = 2850
$ws.Range("I40").Value = 2850
$ws.Range("K40").Value = 2850
$ws.Range("M40").Value = -2714
$ws.Range("H61").Value = 2385.5715
$ws.Range("I61").Value = 2385.5715
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2385.5715
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2385.5715
$ws.Range("I113").Value = 2385.5715
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2385.5715
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 9368.799999999999
$ws.Range("I136").Value = 4864.6665
$ws.Range("K136").Value = 14593.9995
$ws.Range("M136").Value = -12043.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13499.583
$ws.Range("I132").Value = 13998.667
$ws.Range("K132").Value = 41996.001
$ws.Range("M132").Value = -39466.001
